{"js": "// Renumber the routes mentioned in the transcript:\n//   route 5  -> route 4\n//   route 7  -> route 6\n//   route 11 -> route 10\n// (the narrative was updated after the output files were switched to PDF,\n// which re-ordered/renamed the generated routes)\n\nconst body = context.document.body;\n\nasync function replaceText(findText, newText) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n    await context.sync();\n  }\n}\n\nawait replaceText(\n  \"routes 5, 7, and 11\",\n  \"routes 4, 6, and 10\"\n);\nawait replaceText(\"This is route 5.\", \"This is route 4.\");\nawait replaceText(\n  \"This our next suggested route, route 7.\",\n  \"This our next suggested route, route 6.\"\n);\nawait replaceText(\n  \"Finally our last route, number 11.\",\n  \"Finally our last route, number 10.\"\n);\n\n// Word tracks the location of the most recent edit with the hidden\n// \"_GoBack\" bookmark. The first replacement above is the edit the\n// document now points to, so move the bookmark there to match.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst landing = body.search(\"routes 4, 6, and 10\", { matchCase: true });\nlanding.load(\"items\");\nawait context.sync();\nif (landing.items.length > 0) {\n  const endOfEdit = landing.items[0].getRange(\"End\");\n  endOfEdit.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Renumber the routes mentioned in the transcript:\n#   route 5  -> route 4\n#   route 7  -> route 6\n#   route 11 -> route 10\n# (the document's narrative was updated after the output files were\n# switched to PDF, which re-ordered/renamed the generated routes)\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute([ref]$findText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$replaceText, 2) | Out-Null\n}\n\nReplace-Text \"routes 5, 7, and 11\" \"routes 4, 6, and 10\"\nReplace-Text \"This is route 5.\" \"This is route 4.\"\nReplace-Text \"This our next suggested route, route 7.\" \"This our next suggested route, route 6.\"\nReplace-Text \"Finally our last route, number 11.\" \"Finally our last route, number 10.\"\n\n# Word tracks the location of the most recent edit with the hidden\n# \"_GoBack\" bookmark. The first replacement above is the edit the\n# document now points to, so move the bookmark there to match.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$anchor = $d.Content\n$anchor.Find.Execute(\"routes 4, 6, and 10\") | Out-Null\n$goBackRange = $d.Range($anchor.End, $anchor.End)\n$d.Bookmarks.Add(\"_GoBack\", $goBackRange)\n"}
